$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Row 2 and Row 3 data effectively swap places (Maria <-> Emilia), with
# the "Emilia" record (now on row 3) also getting a new username/email
# ("eurbina" -> "eurbina2" @ hotmail instead of gmail).
# ----------------------------------------------------------------------

# Row 2 -> Maria Urbina
$ws.Range("A2").Value = "Maria"
$ws.Range("B2").Value = "Urbina"
$ws.Range("C2").Value = "maurbina"
$ws.Range("D2").Value = "murbina@gmail.com"
$ws.Range("E2").Value = "Biociencias"
$ws.Range("F2").Value = "Parasitología"
$ws.Range("H2").Value = "Bioquímica"
$ws.Range("I2").Value = "Femenino"
$ws.Range("J2").Value = 45869632
$ws.Range("K2").Value = "Licenciado "
$ws.Range("L2").Value = "Investigación celular"
$ws.Range("N2").Value = 4168546321

# Row 3 -> Emilia Urbina (with updated username/email)
$ws.Range("A3").Value = "Emilia"
$ws.Range("B3").Value = "Urbina"
$ws.Range("C3").Value = "eurbina2"
$ws.Range("D3").Value = "eurbina2@hotmail.com"
$ws.Range("E3").Value = "Biociencias"
$ws.Range("F3").Value = "Parasitología"
$ws.Range("H3").Value = "Bioquímica"
$ws.Range("I3").Value = "Femenino"
$ws.Range("J3").Value = 2074967
$ws.Range("K3").Value = "Licenciado"
$ws.Range("L3").Value = "Genética aplicada"
$ws.Range("N3").Value = 4168523651

# ----------------------------------------------------------------------
# Hyperlinks: the mailto targets (rId1 -> eurbina@gmail.com, rId2 ->
# murbina@gmail.com, rId3 -> jmagallanes@gmail.com) stay the same, only
# the displayed text changes for D2/D3. This engine always rebuilds the
# hyperlink collection (in-place edits duplicate entries), so delete all
# and re-add them in the right order with the right display text.
# Adding a hyperlink also overwrites the cell's value with the display
# text, so the real cell values (set above) are re-applied afterwards -
# the diff actually wants D3's cell text ("eurbina2@hotmail.com") to
# differ slightly from the hyperlink's display text ("eurbina2@hotmail").
# ----------------------------------------------------------------------

$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:eurbina@gmail.com", "", "", "murbina@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:murbina@gmail.com", "", "", "eurbina2@hotmail")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:jmagallanes@gmail.com", "", "", "jmagallanes@gmail.com")

# Restore true cell values (D3 differs from its hyperlink display text).
$ws.Range("D2").Value = "murbina@gmail.com"
$ws.Range("D3").Value = "eurbina2@hotmail.com"
$ws.Range("D4").Value = "jmagallanes@gmail.com"

# ----------------------------------------------------------------------
# View state: scroll back to the top-left (A1) and select C3.
# ----------------------------------------------------------------------
$ws.Range("C3").Select()
